$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "P00014"
$ws.Range("B13").Value = "Bague"
$ws.Range("C13").Value = 35000
